$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
